$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-09-02 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-09-03 Wednesday", 2)

$tbl = $d.Tables.Item(1)

$values = @(
    @(1,1,"90÷8=11, 2"),
    @(1,2,"32÷9=3, 5"),
    @(1,3,"24÷7=3, 3"),
    @(1,4,"91÷8=11, 3"),
    @(1,5,"36÷9=4, 0"),
    @(5,1,"61÷3=20, 1"),
    @(5,2,"20÷4=5, 0"),
    @(5,3,"13÷6=2, 1"),
    @(5,4,"32÷5=6, 2"),
    @(5,5,"12÷3=4, 0"),
    @(9,1,"85÷8=10, 5"),
    @(9,2,"76÷4=19, 0"),
    @(9,3,"24÷4=6, 0"),
    @(9,4,"49÷8=6, 1"),
    @(9,5,"85÷9=9, 4"),
    @(13,1,"32÷9=3, 5"),
    @(13,2,"11÷5=2, 1"),
    @(13,3,"39÷4=9, 3"),
    @(13,4,"23÷7=3, 2"),
    @(13,5,"12÷6=2, 0"),
    @(17,1,"38÷6=6, 2"),
    @(17,2,"27÷7=3, 6"),
    @(17,3,"33÷3=11, 0"),
    @(17,4,"41÷3=13, 2"),
    @(17,5,"60÷3=20, 0")
)

foreach ($entry in $values) {
    $row = $entry[0]
    $col = $entry[1]
    $text = $entry[2]
    $cell = $tbl.Cell($row, $col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $text
}
